# Trade #47 closed at 2026-02-17 08:39:14 - unknown UNKNOWN +0.000%
#
# Updates the Summary / Strategy Status roll-up cells for the MarketMaking
# strategy and appends the new closed trade (row 48) to both the
# "All Trades" and "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.56   # Current Capital
$wsSummary.Range("B4").Value = -0.44     # Total P&L $
$wsSummary.Range("B5").Value = -0.19     # Total P&L %
$wsSummary.Range("B6").Value = 47        # Total Trades
$wsSummary.Range("B7").Value = 17        # Winning Trades
$wsSummary.Range("B9").Value = 36.17     # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.56      # Capital
$wsStatus.Range("D4").Value = 47         # Trades
$wsStatus.Range("E4").Value = -0.44      # P&L $
$wsStatus.Range("F4").Value = -0.44      # P&L %
$wsStatus.Range("G4").Value = 36.17      # Win Rate %

# ---------------------------------------------------------------------
# 3) Append the new trade (row 48) to "All Trades" and "MarketMaking"
#    sheets. Copy row 47 down to row 48 first so formatting / cell
#    types stay consistent with the rest of the table, then overwrite
#    the cells that actually changed for this trade.
# ---------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A47:Q47").Copy()
    $ws.Range("A48:Q48").PasteSpecial()

    # Only the cells that actually differ from row 47 need to be
    # overwritten - leaving the rest alone keeps their text/number
    # typing exactly as copied (avoids e.g. re-typing an already
    # correct date string and having it re-interpreted as a date serial).
    $ws.Range("A48").Value = 47
    $ws.Range("C48").Value = "08:39:07"
    $ws.Range("F48").Value = 0.84
    $ws.Range("I48").Value = 2.381
    $ws.Range("J48").Value = 0.02
    $ws.Range("K48").Value = 99.56
}
